# Adds a new "TextBox 13" text box to slide 1 containing the title
# "MSc Research Project - Plan for first month", centred, bold, 20pt,
# matching the OOXML produced when a user draws Insert > Text Box in
# PowerPoint and types the heading (word-wrap off, shape auto-fits to
# the typed text).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Position/size (points) chosen so that, once PowerPoint's autofit-to-text
# recalculates the box, the stored EMU extents land exactly on
# x=4190544 y=1549400 cx=4735142 cy=400110.
$left   = 4190544 / 12700.0
$top    = 1549400 / 12700.0
$width  = 4735142 / 12700.0
$height = 400110  / 12700.0

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)

$tf = $tb.TextFrame
$tf.WordWrap = $false
$tf.AutoSize = 1

$tb.Fill.Visible = $false

$tr = $tf.TextRange
$tr.Text = "MSc Research Project - Plan for first month"
$tr.Font.Size = 20
$tr.Font.Bold = $true
$tr.LanguageID = "de-DE"
$tr.ParagraphFormat.Alignment = 2

# Split into the same runs as the authored slide (one run per spell-checked
# word/space) by touching each sub-range individually.
$runSpans = @(
  @(1, 3),    # "MSc"
  @(4, 25),   # " Research Project - Plan "
  @(29, 3),   # "for"
  @(32, 1),   # " "
  @(33, 5),   # "first"
  @(38, 1),   # " "
  @(39, 5)    # "month"
)
foreach ($span in $runSpans) {
  $run = $tr.Characters($span[0], $span[1])
  $run.Font.Bold = $true
  $run.LanguageID = "de-DE"
}
